# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# The site was regenerated, bumping several interest counters by a small
# amount (the same events appear on both sheets, hence the duplicated set
# of edits).

$wb = $excel.ActiveWorkbook

$updates = @{
    9  = 114
    14 = 449
    20 = 70
    26 = 5969
    31 = 14661
    35 = 90
    36 = 9485
    38 = 4222
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates[$row]
}

$updates4 = @{
    9  = 114
    14 = 449
    21 = 70
    29 = 5969
    34 = 14661
    38 = 90
    39 = 9485
    41 = 4222
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
